$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3 (shifts existing rows 3..22 down to 4..23)
$ws.Rows.Item(3).Insert()

# Copy the formatting of the row-above date label cell onto the new label cell
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Row 3 label (new date bucket 2020-04-01)
$ws.Range("A3").Value = "2020-04-01 00:00:00_diff"

# Row 3 values
$ws.Range("B3").Value = 0.0000000001052939957446597
$ws.Range("C3").Value = 7.869792681105293
$ws.Range("D3").Value = -10.45921331889471
$ws.Range("E3").Value = -2.319131318894706
$ws.Range("F3").Value = -0.5308223188947059
$ws.Range("G3").Value = -3.708941318894706
$ws.Range("H3").Value = -4.004270318894706
